# Updates cryptos list values/percentages per the source diff.
# Some Price-column values are numeric-looking strings (e.g. "0.999") that
# must stay TEXT (matching the original inlineStr cells), so we force the
# cell to a text NumberFormat before assigning, then restore the default
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.966.69'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '3.384.73'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.59'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '3.963.60'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '3.378.28'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '61.054.61'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E18').Value = '  -3.99%  '
$ws.Range('E19').Value = '  -4.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '382.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  -4.14%  '
$ws.Range('D26').Value = '3.523.43'
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.29'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('E32').Value = '  -3.40%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.09'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.01'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').Value = '3.416.21'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0768'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.01%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.457.62'
$ws.Range('E47').Value = '  -5.82%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0266'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.23%  '
